$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 38

$ws.Cells.Item($row, 1).Value = "Linked List Cycle"
$ws.Cells.Item($row, 2).Value = "Linked List"
$ws.Cells.Item($row, 3).Value = "Yes"
$ws.Cells.Item($row, 4).Value = "No"
$ws.Cells.Item($row, 5).Value = "Easy"
$ws.Cells.Item($row, 6).Value = "Easy"
$ws.Cells.Item($row, 7).Value = "141 - Linked List Cycle"

$ws.Hyperlinks.Add($ws.Cells.Item($row, 7), "141 - Linked List Cycle", "", "", "141 - Linked List Cycle")

$ws.Range("A19").Select()
